# Actualización automática hashcode sáb ene  9 03:57:18 CET 2021
# Updates the hashcode values (column B) for the matching keys (column A)
# on the active worksheet, as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "B11";  Old = "f2d8e4b51f987ad30c3ce9202f61284c"; New = "5c22eabfca1000c10d5ae9734c1d5f73" },
    @{ Cell = "B15";  Old = "ff33e05ce4124f5cbf53bf885c4af68d"; New = "88080ae6595d1c7dd8f2eb7d96bd302c" },
    @{ Cell = "B121"; Old = "1f698a3405bd15320c1955fb18b791a3"; New = "af604c1e1543e4f739e1ddec91fc792d" },
    @{ Cell = "B123"; Old = "00bc5663ce8072e226666ca1ef1df7ca"; New = "8f4f6c53fdef9104dcba15b39f4d95d7" },
    @{ Cell = "B223"; Old = "d5ac2fcbfc941b7f2020449d861fcc97"; New = "7409a96f034b294b43ca1eaf002a13de" },
    @{ Cell = "B226"; Old = "938af6cbfeee677639589f82d5bd4183"; New = "ae10dbfd62569e670f39a26b583aeac0" },
    @{ Cell = "B231"; Old = "6cb78729ec7e34e2e3d00f34cefe0a53"; New = "a8277a33cac62b84d037c891f52e8893" },
    @{ Cell = "B248"; Old = "01636a269ec3daa4a66a07c474195c24"; New = "7fe12e55e8ce203f78d3ffbc555bb819" },
    @{ Cell = "B280"; Old = "bb16fd7269b916eb42263675b4c6cc22"; New = "658448c747b484e35c8c628f4f0445f4" },
    @{ Cell = "B282"; Old = "6c51a1a832efc41df8f15ccd452ab72e"; New = "54c07e61a341e99e04fe6b517d06b4a1" },
    @{ Cell = "B417"; Old = "e66f3a6bfc249a167c409d4715943b24"; New = "f7dba8b5804a41c5f250ca2948bc7300" },
    @{ Cell = "B426"; Old = "0841f66eec1f7caf51680bed6f5054c6"; New = "eac6d56063697c4696c84438a0564182" },
    @{ Cell = "B454"; Old = "47a0359a5c7e8641390da1fde64d5f90"; New = "f505ded6eebc9ba8859ce586065452c7" },
    @{ Cell = "B505"; Old = "581b8d3e661eb4b3e04ddfb924f5ec62"; New = "d3592daf27f738b510ee780fb0ef8112" },
    @{ Cell = "B530"; Old = "288cf95f68d978f734aeafd46e059441"; New = "ebdcc369e97c1f7115e22940b60138e4" },
    @{ Cell = "B582"; Old = "e1a6eddc8d35383eb4232b8ff0bfc7db"; New = "010dedb815a41d8bcd8c2a639e19442c" },
    @{ Cell = "B651"; Old = "2cadb5fbe5efbbb1ed6173ec372215f6"; New = "a8e2d2d1f1a82a72bc5bd983c5a57c67" },
    @{ Cell = "B749"; Old = "74c6706958fb2ca6a70191ae5cb0cf45"; New = "45fb08259921ce621ea1b9d74ce91d4b" },
    @{ Cell = "B772"; Old = "f7bc9c59225a602ae8828f104fe09062"; New = "b68d0bd8ec2f9145b55858b5f0c87f7a" },
    @{ Cell = "B778"; Old = "f854d50f95cc0fa3376f704aa1c24f13"; New = "e8eceefdad87a73b082e99b45a82e6d6" },
    @{ Cell = "B803"; Old = "5874b2f19a8d9b4a163aa73e652efcd0"; New = "a614c4c6f98471f7a0538ef78a8c03a2" },
    @{ Cell = "B873"; Old = "370df81783a39107b4484ed3b19cc86d"; New = "53d9fe5445ea312eaed997243ea12d8a" },
    @{ Cell = "B887"; Old = "047c9ef8a9edf23f4edd23b9d7c8476c"; New = "4ac8f3cd6fc25fda71b5c534d99e0243" },
    @{ Cell = "B923"; Old = "ffa24ba1bba8670f108a88c3988afc45"; New = "3aeba3294228492edcacecc052aae677" },
    @{ Cell = "B947"; Old = "6c72477ca0afc0c6cefa7daaf61a33a1"; New = "4eb96f577ee4133c8b6ca757272922da" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.New
}
